$wb = $excel.ActiveWorkbook

# --- Instruction Decoder sheet: fill in newly-documented instructions ---
$ws = $wb.Worksheets.Item("Instruction Decoder")

$ws.Range("I4").Value = "Load constant to Write Addr"

$ws.Range("B5").Value = "LDR"
$ws.Range("I5").Value = "Load from Reg to Write Addr"

$ws.Range("B6").Value = "DISP"
$ws.Range("I6").Value = "Print to TTY"

$ws.Range("B7").Value = "XOR"
$ws.Range("I7").Value = "XOR number and save to Write Addr"

$ws.Range("B8").Value = "AND"
$ws.Range("I8").Value = "AND number and save to Write Addr"

$ws.Range("B9").Value = "OR"
$ws.Range("I9").Value = "OR number and save to Write Addr"

$ws.Range("B10").Value = "ADD"
$ws.Range("I10").Value = "Add number and save to Write Addr"

$ws.Range("B11").Value = "NOT"
$ws.Range("I11").Value = "Inverse bits and save to Write Addr"

$ws.Range("B12").Value = "JUMP"
$ws.Range("I12").Value = "Unconditional jump"

$ws.Range("I13").Value = "Conditional jump"

$ws.Range("B14").Value = "POP"
$ws.Range("I14").Value = "Pop data from stack"

$ws.Range("B16").Value = "SUB"
$ws.Range("I16").Value = "Subtract number and save to Write Addr"

# --- View state: the active sheet/selection moved from "ALU Functions" to
#     "Instruction Decoder", and the selection there is now a single cell. ---
$wsAlu = $wb.Worksheets.Item("ALU Functions")
[void]$wsAlu.Activate()
[void]$wsAlu.Range("C8").Select()

[void]$ws.Activate()
[void]$ws.Range("I4").Select()
